# Auto-generated: update cached market-price / profit figures
# on the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
# Values correspond to refreshed Universalis current-average-price data
# pulled in by the scheduled runner; profit columns (M/N) are recomputed
# to match, and are added/removed depending on whether the HQ/NQ price is > 0.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 50000
$ws.Range("I10").Value = 50000
$ws.Range("K10").Value = 50000
$ws.Range("M10").Value = -49707
$ws.Range("H40").Value = 2338.7334
$ws.Range("I40").Value = 2364.3333
$ws.Range("J40").Value = 2300.3333
$ws.Range("K40").Value = 2364.3333
$ws.Range("L40").Value = 2300.3333
$ws.Range("M40").Value = -2189.3333
$ws.Range("N40").Value = -2650.3333
$ws.Range("H51").Value = 5024.375
$ws.Range("I51").Value = 4200
$ws.Range("J51").Value = 5142.143
$ws.Range("K51").Value = 4200
$ws.Range("L51").Value = 5142.143
$ws.Range("M51").Value = -3716
$ws.Range("N51").Value = -6110.143
$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23498
$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -73488
$ws.Range("H100").Value = 23812182
$ws.Range("I100").Value = 83334340
$ws.Range("J100").Value = 3320
$ws.Range("K100").Value = 83334340
$ws.Range("L100").Value = 3320
$ws.Range("M100").Value = -83333799
$ws.Range("N100").Value = -4402
$ws.Range("H116").Value = 3144.4443
$ws.Range("I116").Value = 3144.4443
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3144.4443
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 297.5556999999999
$ws.Range("N116").ClearContents()
$ws.Range("H138").Value = 9163030
$ws.Range("I138").Value = 1898456
$ws.Range("J138").Value = 16669756
$ws.Range("K138").Value = 5695368
$ws.Range("L138").Value = 50009268
$ws.Range("M138").Value = -5690228
$ws.Range("N138").Value = -50019548

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27078.264
$ws.Range("I2").Value = 27798.756
$ws.Range("J2").Value = 420
$ws.Range("K2").Value = 27798.756
$ws.Range("L2").Value = 420
$ws.Range("M2").Value = -27685.756
$ws.Range("N2").Value = -646
$ws.Range("H36").Value = 16250
$ws.Range("I36").Value = 16250
$ws.Range("K36").Value = 16250
$ws.Range("M36").Value = -15904
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 10064.667
$ws.Range("I63").Value = 11404.571
$ws.Range("J63").Value = 5375
$ws.Range("K63").Value = 11404.571
$ws.Range("L63").Value = 5375
$ws.Range("M63").Value = -10718.571
$ws.Range("N63").Value = -6747
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 10064.667
$ws.Range("I66").Value = 11404.571
$ws.Range("J66").Value = 5375
$ws.Range("K66").Value = 57022.855
$ws.Range("L66").Value = 26875
$ws.Range("M66").Value = -53590.855
$ws.Range("N66").Value = -33739
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H105").Value = 27890
$ws.Range("J105").Value = 27890
$ws.Range("L105").Value = 27890
$ws.Range("N105").Value = -34878
$ws.Range("H116").Value = 27078.264
$ws.Range("I116").Value = 27798.756
$ws.Range("J116").Value = 420
$ws.Range("K116").Value = 27798.756
$ws.Range("L116").Value = 420
$ws.Range("M116").Value = -25504.756
$ws.Range("N116").Value = -5008
$ws.Range("H122").Value = 3635.6667
$ws.Range("I122").Value = 2692
$ws.Range("K122").Value = 8076
$ws.Range("M122").Value = -5626

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27078.264
$ws.Range("I3").Value = 27798.756
$ws.Range("J3").Value = 420
$ws.Range("K3").Value = 27798.756
$ws.Range("L3").Value = 420
$ws.Range("M3").Value = -27684.756
$ws.Range("N3").Value = -648

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5153.787
$ws.Range("J31").Value = 15815.083
$ws.Range("L31").Value = 15815.083
$ws.Range("N31").Value = -16405.083
$ws.Range("H34").Value = 5153.787
$ws.Range("J34").Value = 15815.083
$ws.Range("L34").Value = 15815.083
$ws.Range("N34").Value = -16219.083
$ws.Range("H62").Value = 18723.857
$ws.Range("I62").Value = 26609.889
$ws.Range("J62").Value = 4529
$ws.Range("K62").Value = 26609.889
$ws.Range("L62").Value = 4529
$ws.Range("M62").Value = -25985.889
$ws.Range("N62").Value = -5777
$ws.Range("H65").Value = 18723.857
$ws.Range("I65").Value = 26609.889
$ws.Range("J65").Value = 4529
$ws.Range("K65").Value = 133049.445
$ws.Range("L65").Value = 22645
$ws.Range("M65").Value = -129929.445
$ws.Range("N65").Value = -28885
$ws.Range("H68").Value = 21746.25
$ws.Range("J68").Value = 23995
$ws.Range("L68").Value = 23995
$ws.Range("N68").Value = -25493
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28630
$ws.Range("H71").Value = 21746.25
$ws.Range("J71").Value = 23995
$ws.Range("L71").Value = 71985
$ws.Range("N71").Value = -79473
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -30184
$ws.Range("H106").Value = 33000
$ws.Range("J106").Value = 33000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -35524
$ws.Range("H134").Value = 2998.7742
$ws.Range("I134").Value = 1469.2858
$ws.Range("J134").Value = 6210.7
$ws.Range("K134").Value = 4407.857400000001
$ws.Range("L134").Value = 18632.1
$ws.Range("M134").Value = -1872.857400000001
$ws.Range("N134").Value = -23702.1

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1553.921
$ws.Range("I5").Value = 1082.5714
$ws.Range("J5").Value = 1660.3549
$ws.Range("K5").Value = 3247.7142
$ws.Range("L5").Value = 4981.0647
$ws.Range("M5").Value = -3135.7142
$ws.Range("N5").Value = -5205.0647
$ws.Range("H55").Value = 3816.6667
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3816.6667
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 11450.0001
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -11804.0001
$ws.Range("H60").Value = 1946.7778
$ws.Range("J60").Value = 4244.75
$ws.Range("L60").Value = 12734.25
$ws.Range("N60").Value = -13236.25
$ws.Range("H70").Value = 2922.75
$ws.Range("I70").Value = 1345.5
$ws.Range("K70").Value = 4036.5
$ws.Range("M70").Value = -3721.5
$ws.Range("H73").Value = 2922.75
$ws.Range("I73").Value = 1345.5
$ws.Range("K73").Value = 4036.5
$ws.Range("M73").Value = -2944.5
$ws.Range("H113").Value = 515.375
$ws.Range("I113").Value = 464.33334
$ws.Range("J113").Value = 527.1539
$ws.Range("K113").Value = 1393.00002
$ws.Range("L113").Value = 1581.4617
$ws.Range("M113").Value = 776.9999800000001
$ws.Range("N113").Value = -5921.4617
$ws.Range("H135").Value = 1553.921
$ws.Range("I135").Value = 1082.5714
$ws.Range("J135").Value = 1660.3549
$ws.Range("K135").Value = 9743.142600000001
$ws.Range("L135").Value = 14943.1941
$ws.Range("M135").Value = -7208.142600000001
$ws.Range("N135").Value = -20013.1941

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3274.7
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 3323.1765
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3323.1765
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -3547.1765
$ws.Range("H22").Value = 9300
$ws.Range("J22").Value = 10681.818
$ws.Range("L22").Value = 10681.818
$ws.Range("N22").Value = -11271.818
$ws.Range("H27").Value = 9300
$ws.Range("J27").Value = 10681.818
$ws.Range("L27").Value = 10681.818
$ws.Range("N27").Value = -10895.818
$ws.Range("H40").Value = 3266.3333
$ws.Range("I40").Value = 2722.5
$ws.Range("J40").Value = 3421.7144
$ws.Range("K40").Value = 2722.5
$ws.Range("L40").Value = 3421.7144
$ws.Range("M40").Value = -2586.5
$ws.Range("N40").Value = -3693.7144
$ws.Range("H106").Value = 19790
$ws.Range("J106").Value = 19790
$ws.Range("L106").Value = 19790
$ws.Range("N106").Value = -22314
$ws.Range("H126").Value = 3274.7
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3323.1765
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 9969.529500000001
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -14909.5295

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 500118
$ws.Range("I75").Value = 500118
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 500118
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -499182
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 500118
$ws.Range("I78").Value = 500118
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 1500354
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -1495674
$ws.Range("N78").ClearContents()
$ws.Range("H101").Value = 24933.5
$ws.Range("J101").Value = 24933.5
$ws.Range("L101").Value = 24933.5
$ws.Range("N101").Value = -31423.5
$ws.Range("H103").Value = 355200.66
$ws.Range("J103").Value = 355200.66
$ws.Range("L103").Value = 355200.66
$ws.Range("N103").Value = -357544.66
$ws.Range("H104").Value = 29390
$ws.Range("J104").Value = 29390
$ws.Range("L104").Value = 29390
$ws.Range("N104").Value = -36378
$ws.Range("H105").Value = 20615
$ws.Range("J105").Value = 20615
$ws.Range("L105").Value = 20615
$ws.Range("N105").Value = -27603
$ws.Range("H123").Value = 33333.332
$ws.Range("J123").Value = 33333.332
$ws.Range("L123").Value = 33333.332
$ws.Range("N123").Value = -43133.332
$ws.Range("H126").Value = 78361.62
$ws.Range("I126").Value = 200980.2
$ws.Range("J126").Value = 1725
$ws.Range("K126").Value = 602940.6000000001
$ws.Range("L126").Value = 5175
$ws.Range("M126").Value = -600470.6000000001
$ws.Range("N126").Value = -10115

Write-Output "Applied $(269)  cell updates across 7 worksheets."